# Add Homework Week 1 - populate the "Final Grade" column (F) with the
# average of the Midterm Exam, Midterm Paper, Final Exam and Final Paper
# columns (B:E) for every student row, formatted as a whole number and
# centered - matching the rest of the grade columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column F (Final Grade) currently holds empty, centered placeholder cells.
# Fill it with =AVERAGE(Bn:En) for rows 2 through 25.
$dataRange = $ws.Range("F2:F25")
$dataRange.FormulaR1C1 = "=AVERAGE(RC[-4]:RC[-1])"

# Match the look of the other grade columns: centered, whole-number format.
$dataRange.NumberFormat = "0"
$dataRange.HorizontalAlignment = -4108  # xlCenter

# Leave the cursor on G5, like the finished workbook.
[void]$ws.Range("G5").Select()
